$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "43.711.30"
$ws.Range("E2").Value = "  +0.78%  "

# Row 3
$ws.Range("D3").Value = "2.350.46"
$ws.Range("E3").Value = "  +0.81%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("E5").Value = "  +3.56%  "

# Row 6
Set-TextValue "D6" "235.26"
$ws.Range("E6").Value = "  +1.07%  "

# Row 7
Set-TextValue "D7" "73.68"
$ws.Range("E7").Value = "  +10.70%  "

# Row 8
$ws.Range("E8").Value = "  -0.09%  "

# Row 9
Set-TextValue "D9" "0.541"
$ws.Range("E9").Value = "  +18.97%  "

# Row 10
Set-TextValue "D10" "0.0985"
$ws.Range("E10").Value = "  +1.66%  "

# Row 11
Set-TextValue "D11" "28.33"
$ws.Range("E11").Value = "  +5.32%  "

# Row 12
$ws.Range("E12").Value = "  +1.91%  "

# Row 13
$ws.Range("D13").Value = "2.699.81"
$ws.Range("E13").Value = "  +0.70%  "

# Row 14
Set-TextValue "D14" "16.71"
$ws.Range("E14").Value = "  +7.76%  "

# Row 15
Set-TextValue "D15" "6.69"
$ws.Range("E15").Value = "  +7.37%  "

# Row 16
Set-TextValue "D16" "0.889"
$ws.Range("E16").Value = "  +4.19%  "

# Row 17
$ws.Range("D17").Value = "2.353.29"
$ws.Range("E17").Value = "  +0.70%  "

# Row 18
$ws.Range("D18").Value = "43.741.51"
$ws.Range("E18").Value = "  +1.02%  "

# Row 19
$ws.Range("E19").Value = "  +3.23%  "

# Row 20
Set-TextValue "D20" "77.10"
$ws.Range("E20").Value = "  +4.03%  "

# Row 21
$ws.Range("E21").Value = "  +2.58%  "

# Row 22
Set-TextValue "D22" "253.31"
$ws.Range("E22").Value = "  +1.84%  "

# Row 23
$ws.Range("E23").Value = "  +0.02%  "

# Row 24
$ws.Range("E24").Value = "  -1.62%  "

# Row 25
Set-TextValue "D25" "2.49"
$ws.Range("E25").Value = "  +2.78%  "

# Row 26
Set-TextValue "D26" "10.55"
$ws.Range("E26").Value = "  +6.19%  "

# Row 27
Set-TextValue "D27" "2.25"
$ws.Range("E27").Value = "  -1.32%  "

# Row 28
Set-TextValue "D28" "22.36"
$ws.Range("E28").Value = "  +1.03%  "

# Row 29
$ws.Range("E29").Value = "  +8.53%  "

# Row 30
Set-TextValue "D30" "172.68"
$ws.Range("E30").Value = "  -0.90%  "

# Row 31
Set-TextValue "D31" "0.131"
$ws.Range("E31").Value = "  +1.87%  "

# Row 32
$ws.Range("E32").Value = "  +4.97%  "

# Row 33
$ws.Range("E33").Value = "  +3.06%  "

# Row 34
Set-TextValue "D34" "0.0714"
$ws.Range("E34").Value = "  +3.72%  "

# Row 35
Set-TextValue "D35" "5.18"
$ws.Range("E35").Value = "  +4.30%  "

# Row 36
Set-TextValue "D36" "3.87"
$ws.Range("E36").Value = "  +7.21%  "

# Row 37
$ws.Range("E37").Value = "  -4.26%  "

# Row 38
$ws.Range("E38").Value = "  -1.89%  "

# Row 39
$ws.Range("E39").Value = "  +5.98%  "

# Row 40
Set-TextValue "D40" "19.48"
$ws.Range("E40").Value = "  +6.70%  "

# Row 41
$ws.Range("E41").Value = "  -0.10%  "

# Row 42
Set-TextValue "D42" "8.84"
$ws.Range("E42").Value = "  -2.45%  "

# Row 43
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D43" "1.23"
$ws.Range("E43").Value = "  +1.60%  "

# Row 44
$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D44" "0.0977"
$ws.Range("E44").Value = "  +3.09%  "

# Row 45
$ws.Range("E45").Value = "  -0.85%  "

# Row 46
Set-TextValue "D46" "4.43"
$ws.Range("E46").Value = "  +2.25%  "

# Row 47
$ws.Range("E47").Value = "  +11.69%  "

# Row 48
Set-TextValue "D48" "97.20"
$ws.Range("E48").Value = "  -2.08%  "

# Row 49
$ws.Range("D49").Value = "1.437.14"
$ws.Range("E49").Value = "  -0.73%  "

# Row 50
$ws.Range("E50").Value = "  +1.44%  "

# Row 51
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.573.65"
$ws.Range("E51").Value = "  +0.76%  "
